# Appends 6 new survey rows (203-208) to Sheet1, mirroring the pattern of the
# preceding "no accident occurred" placeholder rows, and updates the active
# selection to reflect where the author left off scrolling/editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 203; Index = 201; Weather = "Day, clear" },
    @{ Row = 204; Index = 202; Weather = "Night, clear" },
    @{ Row = 205; Index = 203; Weather = "Day, clear" },
    @{ Row = 206; Index = 204; Weather = "Day, clear" },
    @{ Row = 207; Index = 205; Weather = "Day, clear" },
    @{ Row = 208; Index = 206; Weather = "Night, clear" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Index            # A - record index
    $ws.Cells.Item($row, 2).Value = "N/A"                # B - Number of vehicles in accident
    $ws.Cells.Item($row, 3).Value = "N/A"                # C - Accident Type
    $ws.Cells.Item($row, 4).Value = "No"                 # D - Person Injury?
    $ws.Cells.Item($row, 5).Value = "No"                 # E - Need for ambulance?
    $ws.Cells.Item($row, 6).Value = "No"                 # F - Need for firetruck?
    $ws.Cells.Item($row, 7).Value = "No"                 # G - Need for Police?
    $ws.Cells.Item($row, 8).Value = "N/A"                # H - Types of vehicles involved
    $ws.Cells.Item($row, 9).Value = "No"                 # I - Fire?
    $ws.Cells.Item($row, 10).Value = $r.Weather          # J - Weather
    $ws.Cells.Item($row, 11).Value = "No"                # K - Low Res/Bad Footage?
    $ws.Cells.Item($row, 12).Value = "no accident occurred"  # L - Other
}

$ws.Range("D185").Select() | Out-Null
